$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CE1 needs the same date-formatted style as CD1 (s="2"); copy format then set value.
$ws.Range("CD1").Copy()
$ws.Range("CE1").PasteSpecial(-4122)
$ws.Range("CE1").Value = 45891

# Append the new CE column values (plain numeric, no explicit style) for rows 4-35.
$ws.Range("CE4").Value = -0.9947480087648897
$ws.Range("CE5").Value = 2.662919374670691
$ws.Range("CE6").Value = 1.584368476884657
$ws.Range("CE7").Value = 1.064474836623286
$ws.Range("CE8").Value = 1.929871341910561
$ws.Range("CE9").Value = 1.886792452830188
$ws.Range("CE10").Value = 1.983439245137664
$ws.Range("CE11").Value = 3.109264853977867
$ws.Range("CE12").Value = 1.76718349407885
$ws.Range("CE13").Value = -0.2039408571514079
$ws.Range("CE14").Value = -0.5379413974455072
$ws.Range("CE15").Value = 0.6949480299733857
$ws.Range("CE16").Value = 1.056232371121646
$ws.Range("CE17").Value = 4.038244551339165
$ws.Range("CE18").Value = 2.996746389634097
$ws.Range("CE19").Value = 0.6207049434715062
$ws.Range("CE20").Value = -5.499559374311513
$ws.Range("CE21").Value = 4.006994025936161
$ws.Range("CE22").Value = 3.833006444382159
$ws.Range("CE23").Value = 0.6530303848022223
$ws.Range("CE24").Value = 0.5093833780160928
$ws.Range("CE25").Value = 2.173913043478271
$ws.Range("CE26").Value = 1.417569507897154
$ws.Range("CE27").Value = 2.136532125205948
$ws.Range("CE28").Value = 3.094914058168263
$ws.Range("CE29").Value = 1.151420329536013
$ws.Range("CE30").Value = 1.01022307078813
$ws.Range("CE31").Value = -4.486182557722218
$ws.Range("CE32").Value = 3.865230460921842
$ws.Range("CE33").Value = 1.886018860188643
$ws.Range("CE34").Value = -0.6770031956444544
$ws.Range("CE35").Value = -0.4695059462808904
